# New Soft Constraint matcher and Output to CSV
# Update the "Test Min" column (M) on Sheet2:
#   - Rows 2-13: apply integer number format and set value to 20
#   - Rows 14-37: set value to 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("M2:M13").NumberFormat = "0"
$ws.Range("M2:M13").Value = 20

$ws.Range("M14:M37").Value = 0

$ws.Range("A1").Select() | Out-Null
